$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Sergipe"
$ws.Cells.Item(2, 2).Value = "Só estuda"
$ws.Cells.Item(2, 3).Value = "31/12/2012"
$ws.Cells.Item(2, 4).Value = 25.60534787983632

$ws.Cells.Item(3, 1).Value = "Sergipe"
$ws.Cells.Item(3, 2).Value = "Só estuda"
$ws.Cells.Item(3, 3).Value = "31/12/2013"
$ws.Cells.Item(3, 4).Value = 26.609715927658

$ws.Cells.Item(4, 1).Value = "Sergipe"
$ws.Cells.Item(4, 2).Value = "Só estuda"
$ws.Cells.Item(4, 3).Value = "31/12/2014"
$ws.Cells.Item(4, 4).Value = 24.87849391650491

$ws.Cells.Item(5, 1).Value = "Sergipe"
$ws.Cells.Item(5, 2).Value = "Só estuda"
$ws.Cells.Item(5, 3).Value = "31/12/2015"
$ws.Cells.Item(5, 4).Value = 25.6

$ws.Cells.Item(6, 1).Value = "Sergipe"
$ws.Cells.Item(6, 2).Value = "Só estuda"
$ws.Cells.Item(6, 3).Value = "31/12/2016"
$ws.Cells.Item(6, 4).Value = 25.8

$ws.Cells.Item(7, 1).Value = "Sergipe"
$ws.Cells.Item(7, 2).Value = "Só estuda"
$ws.Cells.Item(7, 3).Value = "31/12/2017"
$ws.Cells.Item(7, 4).Value = 24.63403429099106

$ws.Cells.Item(8, 1).Value = "Sergipe"
$ws.Cells.Item(8, 2).Value = "Só estuda"
$ws.Cells.Item(8, 3).Value = "31/12/2018"
$ws.Cells.Item(8, 4).Value = 28.40866137529308

$ws.Cells.Item(9, 1).Value = "Sergipe"
$ws.Cells.Item(9, 2).Value = "Só estuda"
$ws.Cells.Item(9, 3).Value = "31/12/2019"
$ws.Cells.Item(9, 4).Value = 26.09824277068888

$ws.Cells.Item(10, 1).Value = "Sergipe"
$ws.Cells.Item(10, 2).Value = "Só estuda"
$ws.Cells.Item(10, 3).Value = "31/12/2020"
$ws.Cells.Item(10, 4).Value = 29.6261672169715

$ws.Cells.Item(11, 1).Value = "Sergipe"
$ws.Cells.Item(11, 2).Value = "Só estuda"
$ws.Cells.Item(11, 3).Value = "31/12/2021"
$ws.Cells.Item(11, 4).Value = 26.50874503070254

$ws.Cells.Item(12, 1).Value = "Sergipe"
$ws.Cells.Item(12, 2).Value = "Só estuda"
$ws.Cells.Item(12, 3).Value = "31/12/2022"
$ws.Cells.Item(12, 4).Value = 28.62124638481935

$ws.Cells.Item(13, 1).Value = "Sergipe"
$ws.Cells.Item(13, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(13, 3).Value = "31/12/2012"
$ws.Cells.Item(13, 4).Value = 11.96439345884951

$ws.Cells.Item(14, 1).Value = "Sergipe"
$ws.Cells.Item(14, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(14, 3).Value = "31/12/2013"
$ws.Cells.Item(14, 4).Value = 12.01231279225067

$ws.Cells.Item(15, 1).Value = "Sergipe"
$ws.Cells.Item(15, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(15, 3).Value = "31/12/2014"
$ws.Cells.Item(15, 4).Value = 11.74810123988309

$ws.Cells.Item(16, 1).Value = "Sergipe"
$ws.Cells.Item(16, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(16, 3).Value = "31/12/2015"
$ws.Cells.Item(16, 4).Value = 10.9

$ws.Cells.Item(17, 1).Value = "Sergipe"
$ws.Cells.Item(17, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(17, 3).Value = "31/12/2016"
$ws.Cells.Item(17, 4).Value = 9.1

$ws.Cells.Item(18, 1).Value = "Sergipe"
$ws.Cells.Item(18, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(18, 3).Value = "31/12/2017"
$ws.Cells.Item(18, 4).Value = 9.678321003268797

$ws.Cells.Item(19, 1).Value = "Sergipe"
$ws.Cells.Item(19, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(19, 3).Value = "31/12/2018"
$ws.Cells.Item(19, 4).Value = 8.01773407447082

$ws.Cells.Item(20, 1).Value = "Sergipe"
$ws.Cells.Item(20, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(20, 3).Value = "31/12/2019"
$ws.Cells.Item(20, 4).Value = 8.97844527129217

$ws.Cells.Item(21, 1).Value = "Sergipe"
$ws.Cells.Item(21, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(21, 3).Value = "31/12/2020"
$ws.Cells.Item(21, 4).Value = 6.617084314546272

$ws.Cells.Item(22, 1).Value = "Sergipe"
$ws.Cells.Item(22, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(22, 3).Value = "31/12/2021"
$ws.Cells.Item(22, 4).Value = 7.825615709242955

$ws.Cells.Item(23, 1).Value = "Sergipe"
$ws.Cells.Item(23, 2).Value = "Estuda e trabalha"
$ws.Cells.Item(23, 3).Value = "31/12/2022"
$ws.Cells.Item(23, 4).Value = 8.845564513098541

$ws.Cells.Item(24, 1).Value = "Sergipe"
$ws.Cells.Item(24, 2).Value = "Só trabalha"
$ws.Cells.Item(24, 3).Value = "31/12/2012"
$ws.Cells.Item(24, 4).Value = 38.2528652224094

$ws.Cells.Item(25, 1).Value = "Sergipe"
$ws.Cells.Item(25, 2).Value = "Só trabalha"
$ws.Cells.Item(25, 3).Value = "31/12/2013"
$ws.Cells.Item(25, 4).Value = 36.90831472075404

$ws.Cells.Item(26, 1).Value = "Sergipe"
$ws.Cells.Item(26, 2).Value = "Só trabalha"
$ws.Cells.Item(26, 3).Value = "31/12/2014"
$ws.Cells.Item(26, 4).Value = 39.77447861099596

$ws.Cells.Item(27, 1).Value = "Sergipe"
$ws.Cells.Item(27, 2).Value = "Só trabalha"
$ws.Cells.Item(27, 3).Value = "31/12/2015"
$ws.Cells.Item(27, 4).Value = 38

$ws.Cells.Item(28, 1).Value = "Sergipe"
$ws.Cells.Item(28, 2).Value = "Só trabalha"
$ws.Cells.Item(28, 3).Value = "31/12/2016"
$ws.Cells.Item(28, 4).Value = 37.4

$ws.Cells.Item(29, 1).Value = "Sergipe"
$ws.Cells.Item(29, 2).Value = "Só trabalha"
$ws.Cells.Item(29, 3).Value = "31/12/2017"
$ws.Cells.Item(29, 4).Value = 35.04438297797804

$ws.Cells.Item(30, 1).Value = "Sergipe"
$ws.Cells.Item(30, 2).Value = "Só trabalha"
$ws.Cells.Item(30, 3).Value = "31/12/2018"
$ws.Cells.Item(30, 4).Value = 32.27002225421842

$ws.Cells.Item(31, 1).Value = "Sergipe"
$ws.Cells.Item(31, 2).Value = "Só trabalha"
$ws.Cells.Item(31, 3).Value = "31/12/2019"
$ws.Cells.Item(31, 4).Value = 37.38420246533146

$ws.Cells.Item(32, 1).Value = "Sergipe"
$ws.Cells.Item(32, 2).Value = "Só trabalha"
$ws.Cells.Item(32, 3).Value = "31/12/2020"
$ws.Cells.Item(32, 4).Value = 31.3529002628708

$ws.Cells.Item(33, 1).Value = "Sergipe"
$ws.Cells.Item(33, 2).Value = "Só trabalha"
$ws.Cells.Item(33, 3).Value = "31/12/2021"
$ws.Cells.Item(33, 4).Value = 32.69833830944476

$ws.Cells.Item(34, 1).Value = "Sergipe"
$ws.Cells.Item(34, 2).Value = "Só trabalha"
$ws.Cells.Item(34, 3).Value = "31/12/2022"
$ws.Cells.Item(34, 4).Value = 35.35236730639497

$ws.Cells.Item(35, 1).Value = "Sergipe"
$ws.Cells.Item(35, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(35, 3).Value = "31/12/2012"
$ws.Cells.Item(35, 4).Value = 24.17739343890454

$ws.Cells.Item(36, 1).Value = "Sergipe"
$ws.Cells.Item(36, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(36, 3).Value = "31/12/2013"
$ws.Cells.Item(36, 4).Value = 24.46965655933733

$ws.Cells.Item(37, 1).Value = "Sergipe"
$ws.Cells.Item(37, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(37, 3).Value = "31/12/2014"
$ws.Cells.Item(37, 4).Value = 23.59892623261636

$ws.Cells.Item(38, 1).Value = "Sergipe"
$ws.Cells.Item(38, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(38, 3).Value = "31/12/2015"
$ws.Cells.Item(38, 4).Value = 25.5

$ws.Cells.Item(39, 1).Value = "Sergipe"
$ws.Cells.Item(39, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(39, 3).Value = "31/12/2016"
$ws.Cells.Item(39, 4).Value = 27.7

$ws.Cells.Item(40, 1).Value = "Sergipe"
$ws.Cells.Item(40, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(40, 3).Value = "31/12/2017"
$ws.Cells.Item(40, 4).Value = 30.6432617277623

$ws.Cells.Item(41, 1).Value = "Sergipe"
$ws.Cells.Item(41, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(41, 3).Value = "31/12/2018"
$ws.Cells.Item(41, 4).Value = 31.30358229601752

$ws.Cells.Item(42, 1).Value = "Sergipe"
$ws.Cells.Item(42, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(42, 3).Value = "31/12/2019"
$ws.Cells.Item(42, 4).Value = 27.53910949268759

$ws.Cells.Item(43, 1).Value = "Sergipe"
$ws.Cells.Item(43, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(43, 3).Value = "31/12/2020"
$ws.Cells.Item(43, 4).Value = 32.40384820561173

$ws.Cells.Item(44, 1).Value = "Sergipe"
$ws.Cells.Item(44, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(44, 3).Value = "31/12/2021"
$ws.Cells.Item(44, 4).Value = 32.96730095060988

$ws.Cells.Item(45, 1).Value = "Sergipe"
$ws.Cells.Item(45, 2).Value = "Não estuda e não trabalha"
$ws.Cells.Item(45, 3).Value = "31/12/2022"
$ws.Cells.Item(45, 4).Value = 27.18082179568723
